$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.905.66"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "2.503.50"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'588.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.72%  "
$ws.Range("D6").Value = "'176.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.19%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +0.90%  "
$ws.Range("E9").Value = "  +3.87%  "
$ws.Range("E10").Value = "  +0.36%  "
$ws.Range("D11").Value = "'0.341"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.01%  "
$ws.Range("D12").Value = "'4.95"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.965.74"
$ws.Range("E13").Value = "  +1.26%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "'25.86"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.29%  "
$ws.Range("D15").Value = "67.684.61"
$ws.Range("E16").Value = "  +1.51%  "
$ws.Range("D17").Value = "2.509.72"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D19").Value = "'7.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.84%  "
$ws.Range("D20").Value = "'351.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("D21").Value = "'4.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.95%  "
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Value = "'70.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E24").Value = "  +1.75%  "
$ws.Range("D25").Value = "'1.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("E26").Value = "  +0.70%  "
$ws.Range("E27").Value = "  +1.41%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").Value = "0.0₃0914"
$ws.Range("E29").Value = "  +1.45%  "
$ws.Range("D30").Value = "'511.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("E31").Value = "  +2.67%  "
$ws.Range("E32").Value = "  +3.23%  "
$ws.Range("D33").Value = "'1.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.10%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("D35").Value = "'0.123"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.88%  "
$ws.Range("D36").Value = "'162.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.11%  "
$ws.Range("E37").Value = "  +1.28%  "
$ws.Range("D38").Value = "'18.69"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("E39").Value = "  +1.49%  "
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("D41").Value = "'1.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.53%  "
$ws.Range("E42").Value = "  +1.20%  "
$ws.Range("D43").Value = "'4.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.28%  "
$ws.Range("D44").Value = "'2.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.00%  "
$ws.Range("D45").Value = "'145.69"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.47%  "
$ws.Range("E46").Value = "  +2.38%  "
$ws.Range("D47").Value = "'0.517"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.91%  "
$ws.Range("E48").Value = "  +2.04%  "
$ws.Range("E49").Value = "  +1.89%  "
$ws.Range("D50").Value = "'0.589"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.06%  "
$ws.Range("E51").Value = "  +0.55%  "
